$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row above row 2 (shifts existing rows 2-20 down to 3-21)
$ws.Rows.Item(2).Insert()

# Populate the new row 2 with the Angular entry
$ws.Range("A2").Value = "Angular"
$ws.Range("B2").Value = "Adopt"
$ws.Range("C2").Value = "Languages & Frameworks"
$ws.Range("D2").Value = $true
$ws.Range("E2").Value = "Angular is a client-based SPA Framework.  I'm listing it here because the Radar draws the rings in order listed in the spreadsheet, wierdly."

$ws.Range("A4").Select()
